# Fixed file names due to https://www.drupal.org/node/2302893
#
# Renames the sample filenames shown in the "flag" routing/menu-links
# diagrams on slides 25, 26, 28 and 29:
#   Flag.menu_links.yml      -> Flag.links.menu.yml
#   Flag.Local_ACTIONs.yml   -> Flag.links.ACTION.yml

$p = $ppt.ActivePresentation

# Slides 25, 26, 28, 29 each have a "TextBox 15" shape that shows the
# "Flag.menu_links.yml" filename.
foreach ($slideIndex in 25, 26, 28, 29) {
    $slide = $p.Slides.Item($slideIndex)
    $shape = $slide.Shapes.Item("TextBox 15")
    $shape.TextFrame.TextRange.Text = "Flag.links.menu.yml"
}

# Slides 28 and 29 additionally have a "TextBox 19" shape that shows the
# "Flag.Local_ACTIONs.yml" filename.
foreach ($slideIndex in 28, 29) {
    $slide = $p.Slides.Item($slideIndex)
    $shape = $slide.Shapes.Item("TextBox 19")
    $shape.TextFrame.TextRange.Text = "Flag.links.ACTION.yml"
}
